# Generate Report for Handback
#
# Before this edit every locale row shows "Ready for handoff" and has no
# target/handback info recorded yet (Latest Target File / Latest Handback
# File are blank, Latest Handback DateTime is the zero date). After running
# the handback, each row's status becomes "Handed back: in sync with en-US",
# and the target file / handback file / handback datetime columns are filled
# in. The Overview sheet mirrors the same status text per locale.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # BGR-encoded RGB(0x64,0x95,0xED) -> matches the workbook's existing HyperLink style

function Set-HandbackCell($ws, $cell, $text, $url) {
    $ws.Range($cell).Value = $text
    $ws.Hyperlinks.Add($ws.Range($cell), $url, "", "", $text) | Out-Null
    $ws.Range($cell).Font.Underline = $true
    $ws.Range($cell).Font.Color = $hyperlinkColor
}

# ---- Overview sheet: refresh the status column for zh-cn / de-de ----
$ov = $wb.Worksheets.Item("Overview")
foreach ($row in 2,3) {
    $ov.Range("B$row").Value = $statusHandedBack
    $ov.Range("C$row").Value = $statusHandedBack
}

# ---- Per-locale detail sheets ----
$sourceMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/01974ad371b0b99c7594f3d5e341418efcd77e59/e2e/a.md"

$locales = @(
    @{
        Sheet      = "zh-cn"
        XlfName    = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        XlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f479930847e626553cc6a8341f68a2982bca9aab/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        Handback   = "2016-02-17 04:18:55"
    },
    @{
        Sheet      = "de-de"
        XlfName    = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        XlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d09dbbb709b2d95f3646f9c968707bae80ea1125/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        Handback   = "2016-02-17 04:19:13"
    }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    foreach ($row in 2,3) {
        # Status (column B) -> handed back
        $ws.Range("B$row").Value = $statusHandedBack

        # Latest Target File (E) and Latest Handback File (F) are now known
        Set-HandbackCell $ws "E$row" "a.md" $sourceMdUrl
        Set-HandbackCell $ws "F$row" $loc.XlfName $loc.XlfUrl

        # Latest Handback DateTime (G)
        $ws.Range("G$row").Value = $loc.Handback
    }
}

Write-Host "Handback report generated."
